$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 7
$ws.Range("H7").Value = 21900
$ws.Range("J7").Value = 21900
$ws.Range("L7").Value = 21900
$ws.Range("N7").Value = -22124

# Row 14
$ws.Range("H14").Value = 21900
$ws.Range("J14").Value = 21900
$ws.Range("L14").Value = 21900
$ws.Range("N14").Value = -22282

# Row 40
$ws.Range("H40").Value = 909.06665
$ws.Range("I40").Value = 899.25
$ws.Range("J40").Value = 912.63635
$ws.Range("K40").Value = 899.25
$ws.Range("L40").Value = 912.63635
$ws.Range("M40").Value = -724.25
$ws.Range("N40").Value = -1262.63635

# Row 74
$ws.Range("H74").Value = 2044500.8
$ws.Range("I74").Value = 2567712.8
$ws.Range("J74").Value = 3973.6
$ws.Range("K74").Value = 2567712.8
$ws.Range("L74").Value = 3973.6
$ws.Range("M74").Value = -2566776.8
$ws.Range("N74").Value = -5845.6

# Row 76
$ws.Range("H76").Value = 38465268
$ws.Range("I76").Value = 71431570
$ws.Range("J76").Value = 4584
$ws.Range("K76").Value = 71431570
$ws.Range("L76").Value = 4584
$ws.Range("M76").Value = -71431255
$ws.Range("N76").Value = -5214

# Row 77
$ws.Range("H77").Value = 2044500.8
$ws.Range("I77").Value = 2567712.8
$ws.Range("J77").Value = 3973.6
$ws.Range("K77").Value = 12838564
$ws.Range("L77").Value = 19868
$ws.Range("M77").Value = -12833884
$ws.Range("N77").Value = -29228

# Row 79
$ws.Range("H79").Value = 38465268
$ws.Range("I79").Value = 71431570
$ws.Range("J79").Value = 4584
$ws.Range("K79").Value = 71431570
$ws.Range("L79").Value = 4584
$ws.Range("M79").Value = -71430478
$ws.Range("N79").Value = -6768

# Row 80
$ws.Range("H80").Value = 9095499
$ws.Range("I80").Value = 3775.5
$ws.Range("K80").Value = 11326.5
$ws.Range("M80").Value = -10328.5

# Row 83
$ws.Range("H83").Value = 9095499
$ws.Range("I83").Value = 3775.5
$ws.Range("K83").Value = 33979.5
$ws.Range("M83").Value = -28987.5

# Row 118
$ws.Range("H118").Value = 1140.45
$ws.Range("I118").Value = 1270
$ws.Range("J118").Value = 1084.9286
$ws.Range("K118").Value = 3810
$ws.Range("L118").Value = 3254.7858
$ws.Range("M118").Value = -2153
$ws.Range("N118").Value = -6568.7858

# Row 137
$ws.Range("H137").Value = 20133.885
$ws.Range("I137").Value = 22575.738
$ws.Range("J137").Value = 1413
$ws.Range("K137").Value = 67727.21400000001
$ws.Range("L137").Value = 4239
$ws.Range("M137").Value = -65177.21400000001
$ws.Range("N137").Value = -9339

# Row 141
$ws.Range("H141").Value = 2331.2163
$ws.Range("I141").Value = 1733.0869
$ws.Range("J141").Value = 3313.8572
$ws.Range("K141").Value = 5199.2607
$ws.Range("L141").Value = 9941.571599999999
$ws.Range("M141").Value = -19.26069999999982
$ws.Range("N141").Value = -20301.5716


$ws = $wb.Worksheets.Item("ARM")

# Row 21
$ws.Range("H21").Value = 24683.5
$ws.Range("I21").Value = 1175
$ws.Range("J21").Value = 36437.75
$ws.Range("K21").Value = 1175
$ws.Range("L21").Value = 36437.75
$ws.Range("M21").Value = -801
$ws.Range("N21").Value = -37185.75

# Row 63
$ws.Range("H63").Value = 2725.625
$ws.Range("I63").Value = 2725.625
$ws.Range("K63").Value = 2725.625
$ws.Range("M63").Value = -2039.625

# Row 66
$ws.Range("H66").Value = 2725.625
$ws.Range("I66").Value = 2725.625
$ws.Range("K66").Value = 13628.125
$ws.Range("M66").Value = -10196.125

# Row 74
$ws.Range("H74").Value = 53163.72
$ws.Range("I74").Value = 62320.816
$ws.Range("J74").Value = 2799.6667
$ws.Range("K74").Value = 62320.816
$ws.Range("L74").Value = 2799.6667
$ws.Range("M74").Value = -61446.816
$ws.Range("N74").Value = -4547.6667

# Row 77
$ws.Range("H77").Value = 53163.72
$ws.Range("I77").Value = 62320.816
$ws.Range("J77").Value = 2799.6667
$ws.Range("K77").Value = 311604.08
$ws.Range("L77").Value = 13998.3335
$ws.Range("M77").Value = -307236.08
$ws.Range("N77").Value = -22734.3335

# Row 132
$ws.Range("H132").Value = 3617660.5
$ws.Range("I132").Value = 4251133
$ws.Range("J132").Value = 1445755.8
$ws.Range("K132").Value = 12753399
$ws.Range("L132").Value = 4337267.4
$ws.Range("M132").Value = -12750869
$ws.Range("N132").Value = -4342327.4


$ws = $wb.Worksheets.Item("BSM")

# Row 32
$ws.Range("H32").Value = 10000
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("N32").Value = -10768

# Row 99
$ws.Range("H99").Value = 802.2727
$ws.Range("I99").Value = 788
$ws.Range("J99").Value = 840.3333
$ws.Range("K99").Value = 788
$ws.Range("L99").Value = 840.3333
$ws.Range("M99").Value = 710
$ws.Range("N99").Value = -3836.3333

# Row 105
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()


$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 31290.572
$ws.Range("I31").Value = 39387.383
$ws.Range("J31").Value = 7899.778
$ws.Range("K31").Value = 39387.383
$ws.Range("L31").Value = 7899.778
$ws.Range("M31").Value = -39092.383
$ws.Range("N31").Value = -8489.778

# Row 34
$ws.Range("H34").Value = 31290.572
$ws.Range("I34").Value = 39387.383
$ws.Range("J34").Value = 7899.778
$ws.Range("K34").Value = 39387.383
$ws.Range("L34").Value = 7899.778
$ws.Range("M34").Value = -39185.383
$ws.Range("N34").Value = -8303.778

# Row 35
$ws.Range("H35").Value = 18500.309
$ws.Range("I35").Value = 1157.7778
$ws.Range("J35").Value = 57521
$ws.Range("K35").Value = 1157.7778
$ws.Range("L35").Value = 57521
$ws.Range("M35").Value = -863.7778000000001
$ws.Range("N35").Value = -58109

# Row 62
$ws.Range("H62").Value = 62502690
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748

# Row 65
$ws.Range("H65").Value = 62502690
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740


$ws = $wb.Worksheets.Item("CUL")

# Row 117
$ws.Range("H117").Value = 709.6
$ws.Range("I117").Value = 550
$ws.Range("J117").Value = 949
$ws.Range("K117").Value = 1650
$ws.Range("L117").Value = 2847
$ws.Range("M117").Value = 1792
$ws.Range("N117").Value = -9731

# Row 129
$ws.Range("H129").Value = 19270682
$ws.Range("I129").Value = 597.7778
$ws.Range("J129").Value = 29472492
$ws.Range("K129").Value = 1793.3334
$ws.Range("L129").Value = 88417476
$ws.Range("M129").Value = 3206.6666
$ws.Range("N129").Value = -88427476


$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 4204.7427
$ws.Range("I70").Value = 4150.276
$ws.Range("K70").Value = 4150.276
$ws.Range("M70").Value = -3880.276

# Row 73
$ws.Range("H73").Value = 4204.7427
$ws.Range("I73").Value = 4150.276
$ws.Range("K73").Value = 4150.276
$ws.Range("M73").Value = -3214.276


$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 4683.8887
$ws.Range("I7").Value = 6712.857
$ws.Range("J7").Value = 3392.7273
$ws.Range("K7").Value = 6712.857
$ws.Range("L7").Value = 3392.7273
$ws.Range("M7").Value = -6600.857
$ws.Range("N7").Value = -3616.7273

# Row 20
$ws.Range("H20").Value = 72503
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 72503
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 72503
$ws.Range("N20").Value = -72955
$ws.Range("M20").ClearContents()

# Row 126
$ws.Range("H126").Value = 4683.8887
$ws.Range("I126").Value = 6712.857
$ws.Range("J126").Value = 3392.7273
$ws.Range("K126").Value = 20138.571
$ws.Range("L126").Value = 10178.1819
$ws.Range("M126").Value = -17668.571
$ws.Range("N126").Value = -15118.1819

